$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Configs")

# --- 1. Insert a new row above the current header row (row 3). ---
# This pushes the existing header + all data rows down by one, and the
# engine automatically slides the existing merged "Year" ranges down too
# (B4:B6 -> B5:B7, B7:B9 -> B8:B10, ... B16:B18 -> B17:B19).
$ws.Rows.Item(3).Insert()

# --- 2. Populate the new "properties" header row. ---
$ws.Range("B3").Value = "properties"
$ws.Range("C3").Value = "origin"
$ws.Range("D3").Value = "Deviation"
$ws.Range("E3").Value = ""

# The inserted row inherits row 2's (title bar) formatting - reset it back
# to a plain, unformatted look (Calibri 11, black, no fill) to match the
# rest of the "Normal 2" bodied table, then add the single right-hand
# border that separates the "Deviation" column from whatever follows.
$hdr = $ws.Range("B3:E3")
$hdr.Interior.Pattern = -4142
$hdr.Borders.Item(7).LineStyle = -4142
$hdr.Borders.Item(8).LineStyle = -4142
$hdr.Borders.Item(9).LineStyle = -4142
$hdr.Borders.Item(10).LineStyle = -4142
$hdr.Font.Name = "Calibri"
$hdr.Font.Size = 11
$hdr.Font.Bold = $false
$hdr.Font.Color = 0

$devCell = $ws.Range("D3")
$devCell.Borders.Item(10).LineStyle = 1
$devCell.Borders.Item(10).Weight = 2
$devCell.Borders.Item(10).ColorIndex = 1

# --- 3. Append a new, blank, unformatted row after the last data row. ---
# Before the insert above, the table's last row was row 19 (2021 / FSST /
# 6.2 / 142800); after the insert it is row 20, so the new trailing blank
# row becomes row 21... but the target layout keeps it immediately below
# the last data row, i.e. row 20 in the *post-insert* numbering actually
# corresponds to inserting a row after the (now) last row 19.
$ws.Rows.Item(20).Insert()
$blank = $ws.Range("B20:E20")
$blank.ClearContents()
$blank.Interior.Pattern = -4142
$blank.Borders.Item(7).LineStyle = -4142
$blank.Borders.Item(8).LineStyle = -4142
$blank.Borders.Item(9).LineStyle = -4142
$blank.Borders.Item(10).LineStyle = -4142
$blank.Font.Name = "Calibri"
$blank.Font.Size = 11
$blank.Font.Bold = $false
$blank.Font.Color = 0
